# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Coliflor" at Feria Lagunitas de
# Puerto Montt. The new record is inserted as row 494, which pushes every
# existing row from 494 down through 605 to 495 through 606 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 494; this shifts rows 494:605 down to 495:606
# and grows the sheet's used range / dimension accordingly (A1:R605 -> A1:R606)
$ws.Rows.Item(494).Insert()

# Populate the newly inserted row with the new record's data
$ws.Cells.Item(494, 1).Value2  = 4
$ws.Cells.Item(494, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(494, 3).Value2  = "Los Lagos"
$ws.Cells.Item(494, 4).Value2  = 45204
$ws.Cells.Item(494, 5).Value2  = 10
$ws.Cells.Item(494, 6).Value2  = 100112008
$ws.Cells.Item(494, 7).Value2  = "Coliflor"
$ws.Cells.Item(494, 8).Value2  = "Sin especificar"
$ws.Cells.Item(494, 9).Value2  = "Primera"
$ws.Cells.Item(494, 10).Value2 = 750
$ws.Cells.Item(494, 11).Value2 = 1500
$ws.Cells.Item(494, 12).Value2 = 1500
$ws.Cells.Item(494, 13).Value2 = 1500
$ws.Cells.Item(494, 14).Value2 = "`$/unidad"
$ws.Cells.Item(494, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(494, 16).Value2 = 1500
$ws.Cells.Item(494, 17).Value2 = 1
$ws.Cells.Item(494, 18).Value2 = "Hortaliza"
